# Commit: "model of machine learning changed"
# The "Catégorie" column (P) on both worksheets holds the machine-learning
# model's predicted category for each row. Re-running / updating the model
# produced new predictions for several rows. The "Multi-usage" category no
# longer appears anywhere in the new predictions, so it disappears from the
# shared-strings table entirely once all referencing cells are updated.

$wb = $excel.ActiveWorkbook

# --- Sheet "dim" ---------------------------------------------------------
$wsDim = $wb.Worksheets.Item("dim")
$wsDim.Range("P2").Value = "CVC"
$wsDim.Range("P4").Value = "CVC"

# --- Sheet "Nuit" ----------------------------------------------------------
$wsNuit = $wb.Worksheets.Item("Nuit")
$wsNuit.Range("P2").Value  = "Eclairage"
$wsNuit.Range("P3").Value  = "Investigation en cours"
$wsNuit.Range("P5").Value  = "Changement de comportement"
$wsNuit.Range("P6").Value  = "Changement de comportement"
$wsNuit.Range("P7").Value  = "CVC"
$wsNuit.Range("P9").Value  = "CVC"
$wsNuit.Range("P10").Value = "Changement de comportement"
$wsNuit.Range("P11").Value = "CVC"
$wsNuit.Range("P12").Value = "Eclairage"
$wsNuit.Range("P13").Value = "Eclairage"
$wsNuit.Range("P14").Value = "Changement de comportement"
$wsNuit.Range("P16").Value = "Changement de comportement"
